$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The dividend/cash deposit row (previously last, buy_date placeholder -693594,
# currency EUR) moves to become row 2. All other (USD) rows shift down by one
# and get their fee amounts slightly adjusted. Forex gains on dividends are now
# tax-free, hence the EUR currency row is pulled to the top and the USD amounts
# are recalculated.

$data = @(
    @(45294, -693594, 2000,    "EUR"),
    @(45293, 44826,   1673.87, "USD"),
    @(45293, 44827,   2324.13, "USD"),
    @(45294, 44827,   932.37,  "USD"),
    @(45294, 44842,   98,      "USD"),
    @(45294, 44867,   330,     "USD"),
    @(45294, 44868,   350,     "USD"),
    @(45294, 44868,   240,     "USD"),
    @(45294, 44879,   47.63,   "USD")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
}
